$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Columns H (Total_Experience) and I (Relevent_Experience) for the
#    existing data rows (2-19) must become text-typed cells holding the
#    same numbers they already hold (e.g. "16" instead of 16).
# ------------------------------------------------------------------
$ws.Range("H2:I19").NumberFormat = "@"
for ($r = 2; $r -le 19; $r++) {
    $h = $ws.Cells.Item($r, 8).Text
    $i = $ws.Cells.Item($r, 9).Text
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
}

# ------------------------------------------------------------------
# 2) Insert two new candidate rows before the current row 20 (which
#    holds candidate 320 / harsh). This pushes that row down to 22.
# ------------------------------------------------------------------
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# New row 20: candidate 318 - cleveland1
# Columns C through M hold text values in this workbook (as produced by
# Excelize originally), so force text formatting before assigning.
$ws.Range("C20:M20").NumberFormat = "@"
$ws.Range("B20").Value = 318
$ws.Range("C20").Value = "2024-01-02"
$ws.Range("D20").Value = "java"
$ws.Range("E20").Value = "cleveland1"
$ws.Range("F20").Value = "55667788"
$ws.Range("G20").Value = "cleveland1@gmail.com"
$ws.Range("H20").Value = "16"
$ws.Range("I20").Value = "15"
$ws.Range("J20").Value = "nasa corporation"
$ws.Range("K20").Value = "3"
$ws.Range("L20").Value = "upgraded for further interview level"
$ws.Range("M20").Value = "shortlisted"

# New row 21: candidate 319 - surmak
$ws.Range("C21:M21").NumberFormat = "@"
$ws.Range("B21").Value = 319
$ws.Range("C21").Value = "2024-01-02"
$ws.Range("D21").Value = "java"
$ws.Range("E21").Value = "surmak"
$ws.Range("F21").Value = "55667788"
$ws.Range("G21").Value = "surmka@gmail.com"
$ws.Range("H21").Value = "16"
$ws.Range("I21").Value = "15"
$ws.Range("J21").Value = "umbrala corporation"
$ws.Range("K21").Value = "3"
$ws.Range("L21").Value = "upgraded for further interview level"
$ws.Range("M21").Value = "shortlisted"

# ------------------------------------------------------------------
# 3) The row that shifted down to 22 (candidate 320 / harsh) also needs
#    its H/I columns converted to text, same value as before.
# ------------------------------------------------------------------
$ws.Range("H22:I22").NumberFormat = "@"
$h22 = $ws.Cells.Item(22, 8).Text
$i22 = $ws.Cells.Item(22, 9).Text
$ws.Cells.Item(22, 8).Value = $h22
$ws.Cells.Item(22, 9).Value = $i22
